$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, matching style of existing header cells (e.g. E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("F1").Value = "time_taken"

# Add the time_taken values for each data row, as text
$timestamps = @(
    "2021-10-05 10:51:23.379086",
    "2021-10-05 10:51:23.379095",
    "2021-10-05 10:51:23.379099",
    "2021-10-05 10:51:23.379101",
    "2021-10-05 10:51:23.379104",
    "2021-10-05 10:51:23.379107",
    "2021-10-05 10:51:23.379109",
    "2021-10-05 10:51:23.379112"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
